$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-enter the "NVARCHAR2" label in D15 (DbType column). The text itself is
# unchanged, but re-committing it causes the workbook's shared-string table to
# be rebuilt/deduped against the other existing "NVARCHAR2" entries (D25/D39),
# dropping the now-redundant duplicate - matching the sharedStrings.xml
# uniqueCount drop (144 -> 143) seen in the authored diff.
$ws.Range("D15").Value = "NVARCHAR2"

# Move the live selection/cursor to D23, matching the sheetView's recorded
# selection after the edit (previously G28, with the view scrolled to A28).
$ws.Range("D23").Select() | Out-Null
